$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of A16 (bold/centered/bordered) onto the newly added A17:A19 cells
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.068666319289838
$ws.Range("D10").Value = 1.067777675235372
$ws.Range("E10").Value = 0.9487134065268789
$ws.Range("F10").Value = 1.068666319289838
$ws.Range("G10").Value = 1.010109180398649
$ws.Range("H10").Value = 0.9379565520909054
$ws.Range("I10").Value = 0.981474051889873
$ws.Range("J10").Value = 1.067777675235372
$ws.Range("K10").Value = 1.008245540881125
$ws.Range("L10").Value = 1.038455930085482
$ws.Range("M10").Value = 1.002449530905253

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.8727383293211178
$ws.Range("D11").Value = 0.9503777337269445
$ws.Range("E11").Value = 1.092976826136958
$ws.Range("F11").Value = 0.8727383293211178
$ws.Range("G11").Value = 0.8663765038449267
$ws.Range("H11").Value = 1.3892266471813
$ws.Range("I11").Value = 1.030183596349878
$ws.Range("J11").Value = 0.9503777337269445
$ws.Range("K11").Value = 1.021677279931951
$ws.Range("L11").Value = 0.9472078046265346
$ws.Range("M11").Value = 1.033646606093521

# Row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.8721444005103356
$ws.Range("D12").Value = 0.952504810201553
$ws.Range("E12").Value = 1.092481039671854
$ws.Range("F12").Value = 0.8721444005103356
$ws.Range("G12").Value = 0.8676400895155442
$ws.Range("H12").Value = 1.387701161620881
$ws.Range("I12").Value = 1.029552540133116
$ws.Range("J12").Value = 0.952504810201553
$ws.Range("K12").Value = 1.022492924936703
$ws.Range("L12").Value = 0.9473186627235195
$ws.Range("M12").Value = 1.033670673608881

# Row 13
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.8725508905270002
$ws.Range("D13").Value = 0.9509374295599213
$ws.Range("E13").Value = 1.092850609788978
$ws.Range("F13").Value = 0.8725508905270002
$ws.Range("G13").Value = 0.866683383265007
$ws.Range("H13").Value = 1.388855318751441
$ws.Range("I13").Value = 1.02995674526145
$ws.Range("J13").Value = 0.9509374295599213
$ws.Range("K13").Value = 1.02189401967445
$ws.Range("L13").Value = 0.9472224551007249
$ws.Range("M13").Value = 1.033639062858966

# Row 14
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.6902359999999996
$ws.Range("D14").Value = 1.340892
$ws.Range("E14").Value = 1.044004000000001
$ws.Range("F14").Value = 0.6902359999999996
$ws.Range("G14").Value = 1.068255999999998
$ws.Range("H14").Value = 1.243235999999999
$ws.Range("I14").Value = 0.9573399999999992
$ws.Range("J14").Value = 1.340892
$ws.Range("K14").Value = 1.192448
$ws.Range("L14").Value = 0.941342
$ws.Range("M14").Value = 1.057327333333333

# Row 15
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.61
$ws.Range("D15").Value = 1.408349999999997
$ws.Range("E15").Value = 1.0582125
$ws.Range("F15").Value = 0.61
$ws.Range("G15").Value = 1.1
$ws.Range("H15").Value = 1.28
$ws.Range("I15").Value = 0.96
$ws.Range("J15").Value = 1.408349999999997
$ws.Range("K15").Value = 1.233281249999998
$ws.Range("L15").Value = 0.9216406249999992
$ws.Range("M15").Value = 1.069427083333333

# Row 16
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.7740280868864
$ws.Range("D16").Value = 1.233587690393598
$ws.Range("E16").Value = 1.030545201049596
$ws.Range("F16").Value = 0.7740280868864
$ws.Range("G16").Value = 1.055304665292798
$ws.Range("H16").Value = 1.155482233753597
$ws.Range("I16").Value = 0.9733050576895997
$ws.Range("J16").Value = 1.233587690393598
$ws.Range("K16").Value = 1.132066445721597
$ws.Range("L16").Value = 0.9530472663039986
$ws.Range("M16").Value = 1.037042155844265

# Row 17
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9935122917394891
$ws.Range("D17").Value = 0.9956782041916629
$ws.Range("E17").Value = 0.9958302196493651
$ws.Range("F17").Value = 0.9935122917394891
$ws.Range("G17").Value = 0.9951069915089162
$ws.Range("H17").Value = 0.9958445734285467
$ws.Range("I17").Value = 0.9948297475606834
$ws.Range("J17").Value = 0.9956782041916629
$ws.Range("K17").Value = 0.995754211920514
$ws.Range("L17").Value = 0.9946332518300015
$ws.Range("M17").Value = 0.9951336713464439

# Row 18
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.011298874880657
$ws.Range("D18").Value = 0.9881835490611545
$ws.Range("E18").Value = 0.9900460804314046
$ws.Range("F18").Value = 1.011298874880657
$ws.Range("G18").Value = 0.9909868733543298
$ws.Range("H18").Value = 0.9894811700845939
$ws.Range("I18").Value = 0.996533177215208
$ws.Range("J18").Value = 0.9881835490611545
$ws.Range("K18").Value = 0.9891148147462796
$ws.Range("L18").Value = 1.000206844813468
$ws.Range("M18").Value = 0.9944216208378913

# Row 19
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 1.00543197511227
$ws.Range("D19").Value = 0.9755622242192407
$ws.Range("E19").Value = 0.997120158705321
$ws.Range("F19").Value = 1.00543197511227
$ws.Range("G19").Value = 0.9844891689177391
$ws.Range("H19").Value = 0.9971856438074491
$ws.Range("I19").Value = 0.9994419832810186
$ws.Range("J19").Value = 0.9755622242192407
$ws.Range("K19").Value = 0.9863411914622808
$ws.Range("L19").Value = 0.9958865832872755
$ws.Range("M19").Value = 0.9932051923405064
